# Automatische test-sync: 2025-06-24 21:26:50
# Adds a new log row (30) to the "Logs" sheet and updates the
# "Retour / Terugbetaling" tally on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

$newRow = 30

$ws.Range("A$newRow").Value = "Verkeerde maat"
$ws.Range("B$newRow").Value = "mailmind.test@zohomail.eu"
$ws.Range("C$newRow").Value = "Hoi, ik heb een product ontvangen maar de maat is verkeerd. Mag ik deze omruilen?`nSent using {0}"
$ws.Range("D$newRow").Value = "Retour / Terugbetaling"
$ws.Range("E$newRow").Value = "Beste klant,`nBedankt voor je bericht. Het spijt me te horen dat de maat van het ontvangen product niet klopt. Om je te helpen met de omruiling, hebben we wat extra informatie nodig. Zou je alsjeblieft je bestelnummer en de juiste maat die je wilt ontvangen kunnen doorgeven? Zodra we deze gegevens hebben, zullen we de omruiling in gang zetten.`nMet vriendelijke groet,`n[Bedrijfsnaam]"
$ws.Range("F$newRow").Value = "2025-06-24 21:26:13"
$ws.Range("G$newRow").Value = "Ja"

# Writing multi-line values auto-sizes the new row's height; AutoFit()
# clears the stale explicit height again so the row stays at the sheet
# default (matching every other row in the log).
$ws.Rows.Item($newRow).AutoFit()

# Extend the two conditional-formatting blocks (Categorie / Beantwoord
# columns) so they keep covering the whole data range, D2:D30 / G2:G30.
$catCF = $ws.Range("D2:D29").FormatConditions.Item(1)
$catCF.ModifyAppliesToRange($ws.Range("D2:D$newRow"))

$answeredCF = $ws.Range("G2:G29").FormatConditions.Item(1)
$answeredCF.ModifyAppliesToRange($ws.Range("G2:G$newRow"))

# Bump the "Retour / Terugbetaling" count on the Dashboard sheet.
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 13
